$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shuffle the practice-round rows (2-17): each target row is populated with
# the problem/response/answer that previously lived at a different row.
# Mapping is new_row -> source_row (from the original, pre-shuffle layout).
$rowMap = @{
    2  = 4
    3  = 13
    4  = 10
    5  = 5
    6  = 15
    7  = 7
    8  = 8
    9  = 9
    10 = 14
    11 = 11
    12 = 12
    13 = 3
    14 = 2
    15 = 6
    16 = 16
    17 = 17
}

# Snapshot original values (A: problem text, B: correct-response bool, C: suggested answer)
$orig = @{}
foreach ($r in 2..17) {
    $orig[$r] = @(
        $ws.Cells.Item($r, 1).Value(),
        $ws.Cells.Item($r, 2).Value(),
        $ws.Cells.Item($r, 3).Value()
    )
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $orig[$srcRow]
    $ws.Cells.Item($destRow, 1).Value = $vals[0]
    $ws.Cells.Item($destRow, 2).Value = $vals[1]
    $ws.Cells.Item($destRow, 3).Value = $vals[2]
}

# Match the final selection left in the sheet: whole row 13 selected, active cell A13.
$ws.Range("A13:XFD13").Select()
